$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "42.569.21"
$ws.Range("E2").Value = "  -1.10%  "

# Row 3
$ws.Range("D3").Value = "2.288.05"
$ws.Range("E3").Value = "  -0.41%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
Set-TextValue $ws.Range("D5") "305.93"
$ws.Range("E5").Value = "  +1.91%  "

# Row 6
Set-TextValue $ws.Range("D6") "95.76"
$ws.Range("E6").Value = "  -2.27%  "

# Row 7
$ws.Range("E7").Value = "  -3.13%  "

# Row 8
$ws.Range("E8").Value = "  +0.03%  "

# Row 10
Set-TextValue $ws.Range("D10") "35.08"
$ws.Range("E10").Value = "  -2.97%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.0783"
$ws.Range("E11").Value = "  -0.65%  "

# Row 12
Set-TextValue $ws.Range("D12") "18.23"
$ws.Range("E12").Value = "  +3.07%  "

# Row 13
$ws.Range("E13").Value = "  +1.00%  "

# Row 14
Set-TextValue $ws.Range("D14") "6.69"
$ws.Range("E14").Value = "  -2.37%  "

# Row 15
$ws.Range("D15").Value = "2.643.86"
$ws.Range("E15").Value = "  -0.42%  "

# Row 16
$ws.Range("D16").Value = "2.284.59"
$ws.Range("E16").Value = "  -0.54%  "

# Row 17
$ws.Range("E17").Value = "  -1.41%  "

# Row 18
$ws.Range("D18").Value = "42.479.72"
$ws.Range("E18").Value = "  -1.03%  "

# Row 19
Set-TextValue $ws.Range("D19") "12.84"
$ws.Range("E19").Value = "  +0.54%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0892"
$ws.Range("E20").Value = "  -2.40%  "

# Row 21
Set-TextValue $ws.Range("D21") "6.00"
$ws.Range("E21").Value = "  -1.82%  "

# Row 22
Set-TextValue $ws.Range("D22") "66.85"
$ws.Range("E22").Value = "  -3.15%  "

# Row 23
Set-TextValue $ws.Range("D23") "235.15"
$ws.Range("E23").Value = "  -0.92%  "

# Row 24
$ws.Range("E24").Value = "  -0.51%  "

# Row 25
$ws.Range("E25").Value = "  +0.16%  "

# Row 26
Set-TextValue $ws.Range("D26") "2.44"
$ws.Range("E26").Value = "  +0.68%  "

# Row 27
$ws.Range("E27").Value = "  +0.15%  "

# Row 28
Set-TextValue $ws.Range("D28") "24.94"
$ws.Range("E28").Value = "  -0.14%  "

# Row 29
$ws.Range("E29").Value = "  +6.65%  "

# Row 30
Set-TextValue $ws.Range("D30") "165.82"
$ws.Range("E30").Value = "  +0.32%  "

# Row 31
$ws.Range("E31").Value = "  -1.17%  "

# Row 32
Set-TextValue $ws.Range("D32") "32.52"
$ws.Range("E32").Value = "  -1.75%  "

# Row 33
Set-TextValue $ws.Range("D33") "1.00"
$ws.Range("E33").Value = "  +0.11%  "

# Row 34
Set-TextValue $ws.Range("D34") "4.68"
$ws.Range("E34").Value = "  -1.20%  "

# Row 35
Set-TextValue $ws.Range("D35") "4.94"
$ws.Range("E35").Value = "  -2.63%  "

# Row 36
Set-TextValue $ws.Range("D36") "17.43"
$ws.Range("E36").Value = "  -2.64%  "

# Row 37
Set-TextValue $ws.Range("D37") "2.39"
$ws.Range("E37").Value = "  -0.86%  "

# Row 38
Set-TextValue $ws.Range("D38") "0.0687"
$ws.Range("E38").Value = "  -1.04%  "

# Row 39
$ws.Range("E39").Value = "  -0.91%  "

# Row 40
Set-TextValue $ws.Range("D40") "1.73"
$ws.Range("E40").Value = "  -2.20%  "

# Row 41
$ws.Range("E41").Value = "  -1.91%  "

# Row 42
$ws.Range("E42").Value = "  -3.39%  "

# Row 43
$ws.Range("D43").Value = "1.990.26"
$ws.Range("E43").Value = "  -0.95%  "

# Row 44
Set-TextValue $ws.Range("D44") "0.0277"
$ws.Range("E44").Value = "  -2.84%  "

# Row 45
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D45") "9.97"
$ws.Range("E45").Value = "  -2.95%  "

# Row 46
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D46") "17.93"
$ws.Range("E46").Value = "  +2.68%  "

# Row 47
$ws.Range("E47").Value = "  -10.33%  "

# Row 48
Set-TextValue $ws.Range("D48") "2.76"
$ws.Range("E48").Value = "  -2.21%  "

# Row 49
$ws.Range("E49").Value = "  +8.51%  "

# Row 50
Set-TextValue $ws.Range("D50") "53.57"
$ws.Range("E50").Value = "  -0.90%  "

# Row 51
$ws.Range("D51").Value = "2.510.58"
$ws.Range("E51").Value = "  -0.51%  "
